{"js": "// Facilitators guidelines - Moebius.docx: English -> Swahili (Kenya) translation\n// 1) Translate the visible English labels/phrases to Swahili.\n// 2) Flip the document's default proofing language from Swahili (Tanzania)\n//    to Swahili (Kenya).\n\nconst body = context.document.body;\n\n// Ordered (longest-search-string-first where one string is a substring of\n// another, e.g. \"General VMC Video Introduction\" contains \"Video\n// Introduction\") list of English -> Swahili replacements. Each left-hand\n// side is the exact, unique text of a table-cell run in the document, so a\n// plain search + replace is safe and won't touch anything else.\nconst pairs = [\n  [\"Video Title\", \"Kichwa cha Video\"],\n  [\"Topic\", \"Mada\"],\n  [\"Geometry\", \"Jiometri\"],\n  [\"Aim(s)\", \"Malengo\"],\n  [\"Length\", \"Urefu\"],\n  [\"Camp Location\", \"Mahali pa Kambi\"],\n  [\"Facilitators\", \"Wawezeshaji\"],\n  [\"N. of students\", \"N. ya wanafunzi\"],\n  [\"Date\", \"Tarehe\"],\n  [\"Resources\", \"Rasilimali\"],\n  [\"needed\", \"inahitajika\"],\n  [\"Preparations\", \"Maandalizi\"],\n  [\"Video time\", \"Muda wa video\"],\n  [\"What facilitator does\", \"Mwezeshaji anafanya nini\"],\n  [\"What learners do\", \"Wanachofanya wanafunzi\"],\n  [\"General VMC Video Introduction\", \"Utangulizi Mkuu wa Video ya VMC\"],\n  [\"Video Introduction\", \"Utangulizi wa Video\"],\n  [\"Introduction of the first experiment\", \"Utangulizi wa jaribio la kwanza\"],\n  [\"Assist the process, provoke thoughts\", \"Kusaidia mchakato, kuchochea mawazo\"],\n];\n\nfor (const [src, target] of pairs) {\n  const results = body.search(src, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(target, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Document default language: sw-TZ (Swahili, Tanzania) -> sw-KE (Swahili, Kenya).\nconst styles = context.document.getStyles();\nconst normalStyle = styles.getByNameOrNullObject(\"Normal\");\nnormalStyle.load(\"nameLocal\");\nawait context.sync();\nnormalStyle.languageId = \"sw-KE\";\nawait context.sync();\n", "ps1": "# Facilitators guidelines - Moebius.docx: English -> Swahili (Kenya) translation\n# 1) Translate the visible English labels/phrases to Swahili.\n# 2) Flip the document's default proofing language from Swahili (Tanzania)\n#    to Swahili (Kenya).\n\n$d = $word.ActiveDocument\n\n# Ordered (longest-search-string-first where one string is a substring of\n# another, e.g. \"General VMC Video Introduction\" contains \"Video\n# Introduction\") list of English -> Swahili replacements. Each left-hand\n# side is the exact, unique text of a table-cell run in the document, so a\n# simple Find/Replace-all is safe and won't touch anything else.\n$pairs = [ordered]@{\n  \"Video Title\"                           = \"Kichwa cha Video\"\n  \"Topic\"                                  = \"Mada\"\n  \"Geometry\"                               = \"Jiometri\"\n  \"Aim(s)\"                                 = \"Malengo\"\n  \"Length\"                                 = \"Urefu\"\n  \"Camp Location\"                          = \"Mahali pa Kambi\"\n  \"Facilitators\"                           = \"Wawezeshaji\"\n  \"N. of students\"                         = \"N. ya wanafunzi\"\n  \"Date\"                                    = \"Tarehe\"\n  \"Resources\"                              = \"Rasilimali\"\n  \"needed\"                                  = \"inahitajika\"\n  \"Preparations\"                           = \"Maandalizi\"\n  \"Video time\"                             = \"Muda wa video\"\n  \"What facilitator does\"                  = \"Mwezeshaji anafanya nini\"\n  \"What learners do\"                       = \"Wanachofanya wanafunzi\"\n  \"General VMC Video Introduction\"         = \"Utangulizi Mkuu wa Video ya VMC\"\n  \"Video Introduction\"                     = \"Utangulizi wa Video\"\n  \"Introduction of the first experiment\"   = \"Utangulizi wa jaribio la kwanza\"\n  \"Assist the process, provoke thoughts\"   = \"Kusaidia mchakato, kuchochea mawazo\"\n}\n\nforeach ($src in $pairs.Keys) {\n  $target = $pairs[$src]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($src, $false, $false, $false, $false, $false, $true, 1, $false, $target, 2)\n}\n\n# Document default language: sw-TZ (Swahili, Tanzania) -> sw-KE (Swahili, Kenya).\n$d.Styles(\"Normal\").LanguageID = \"sw-KE\"\n"}
